$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 85, shifting existing rows 85-88 down to 86-89
# (row 85's old blank D/E/F placeholder, the sums and the labels all move
# down by one; the inserted row inherits formatting from row 84 above it).
$ws.Rows.Item(85).Insert()

$ws.Range("A85").Value = 2014
$ws.Range("B85").Value = 3
$ws.Range("C85").Value = 21
$ws.Range("D85").Value = 0.63888888888888895
$ws.Range("E85").Value = 0.75
$ws.Range("F85").Formula = "=(E85-D85)*24*60"
$ws.Range("G85").Formula = "=F85/60"

# update selection to match diff
$ws.Range("A86").Select()
